{"js": "// Apply the \"Peas Pods Drying\" edits:\n//  1. \"Pea Pods\" line: replace the leading tab with 16 spaces.\n//  2. Several instruction paragraphs that were split across multiple\n//     runs (because of earlier in-place word edits) get their runs\n//     collapsed back down to a single run holding the full sentence\n//     text (formatting-neutral re-typing / cleanup).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Map of paragraph index -> full replacement text (same text already\n// present, just re-issued as one run instead of several).\nconst replacements = {\n  1: \"                Pea Pods\\\\n\",\n  4: \"Before blanching, prepare an ice water bath to quickly cool the peas and stop the cooking process.\\\\n\",\n  6: \"Blanch pea pods in boiling water for 1 minute.\\\\n\",\n  10: \"Once the time is up, immediately immersing the peas in an ice bath until cool.\\\\n\",\n  16: \"In a dehydrator, stir the peas around after about 2 hours to allow for even drying.\\\\n\",\n  18: \"In an oven, it is necessary to stir the peas several times during drying to prevent uneven drying.\\\\n\",\n};\n\nfor (const key of Object.keys(replacements)) {\n  const idx = parseInt(key, 10);\n  items[idx].insertText(replacements[key], \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Apply the \"Peas Pods Drying\" edits:\n#  1. \"Pea Pods\" line: replace the leading tab with 16 spaces.\n#  2. Several instruction paragraphs that were split across multiple\n#     runs (because of earlier in-place word edits) get their runs\n#     collapsed back down to a single run holding the full sentence\n#     text (formatting-neutral re-typing / cleanup).\n$d = $word.ActiveDocument\n\n$replacements = @{\n    2  = \"                Pea Pods\\n\"\n    5  = \"Before blanching, prepare an ice water bath to quickly cool the peas and stop the cooking process.\\n\"\n    7  = \"Blanch pea pods in boiling water for 1 minute.\\n\"\n    11 = \"Once the time is up, immediately immersing the peas in an ice bath until cool.\\n\"\n    17 = \"In a dehydrator, stir the peas around after about 2 hours to allow for even drying.\\n\"\n    19 = \"In an oven, it is necessary to stir the peas several times during drying to prevent uneven drying.\\n\"\n}\n\nforeach ($idx in $replacements.Keys) {\n    $pRange = $d.Paragraphs($idx).Range\n    $fullRange = $d.Range($pRange.Start, $pRange.End)\n    $fullRange.Text = $replacements[$idx]\n}\n"}
